$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '25.791.31'; E = '  -1.68%  ' }
    @{ Row = 3; D = '1.615.94'; E = '  -3.50%  ' }
    @{ Row = 4; D = '1.001'; E = '  -0.22%  ' }
    @{ Row = 5; D = '208.65'; E = '  -1.50%  ' }
    @{ Row = 6; D = '0.5188'; E = '  -1.48%  ' }
    @{ Row = 7; D = '1.001'; E = '  -0.20%  ' }
    @{ Row = 8; D = '0.2563'; E = '  -3.40%  ' }
    @{ Row = 9; D = '0.06247'; E = '  -0.50%  ' }
    @{ Row = 10; D = '20.27'; E = '  -4.92%  ' }
    @{ Row = 11; D = '0.07527'; E = '  -0.54%  ' }
    @{ Row = 12; D = '1.613.19'; E = '  -3.73%  ' }
    @{ Row = 13; D = '4.354'; E = '  -2.46%  ' }
    @{ Row = 14; D = '1.834.68'; E = '  -3.63%  ' }
    @{ Row = 15; D = '0.5421'; E = '  -3.37%  ' }
    @{ Row = 16; D = '0.0₅7878'; E = '  -1.59%  ' }
    @{ Row = 17; D = '63.83'; E = '  -4.68%  ' }
    @{ Row = 18; D = '25.785.72'; E = '  -0.97%  ' }
    @{ Row = 19; D = '1.001'; E = '  -0.14%  ' }
    @{ Row = 20; D = '4.625'; E = '  -3.98%  ' }
    @{ Row = 21; D = '183.27'; E = '  -2.43%  ' }
    @{ Row = 22; D = '10.01'; E = '  -3.97%  ' }
    @{ Row = 23; D = '6.042'; E = '  -2.74%  ' }
    @{ Row = 24; D = '1.002'; E = '  -0.18%  ' }
    @{ Row = 25; D = '144.30'; E = '  -3.68%  ' }
    @{ Row = 26; D = '0.1197'; E = '  -4.48%  ' }
    @{ Row = 27; D = '7.341'; E = '  -3.21%  ' }
    @{ Row = 28; D = '15.44'; E = '  -3.24%  ' }
    @{ Row = 29; D = '1.354'; E = '  -0.45%  ' }
    @{ Row = 30; D = '0.05846'; E = '  -6.06%  ' }
    @{ Row = 31; D = '1.235'; E = '  -3.86%  ' }
    @{ Row = 32; D = '3.364'; E = '  -4.07%  ' }
    @{ Row = 33; D = '3.334'; E = '  -2.83%  ' }
    @{ Row = 34; D = '1.593'; E = '  -2.51%  ' }
    @{ Row = 35; D = '0.9651'; E = '  -3.81%  ' }
    @{ Row = 36; D = '2.380'; E = '  -1.30%  ' }
    @{ Row = 37; D = '2.709'; E = '  -1.66%  ' }
    @{ Row = 38; D = '0.5721'; E = '  -5.57%  ' }
    @{ Row = 39; D = '0.01580'; E = '  -2.43%  ' }
    @{ Row = 40; D = '1.002'; E = '  -0.47%  ' }
    @{ Row = 41; D = '0.8391'; E = '  -3.67%  ' }
    @{ Row = 42; D = '5.661'; E = '  -7.43%  ' }
    @{ Row = 43; D = '1.014.17'; E = '  -7.96%  ' }
    @{ Row = 44; D = '99.17'; E = '  -0.69%  ' }
    @{ Row = 45; D = '1.761.24'; E = '  -3.52%  ' }
    @{ Row = 46; D = '0.0₈109'; E = '  -2.04%  ' }
    @{ Row = 47; D = '1.004'; E = '  -0.22%  ' }
    @{ Row = 48; D = '54.15'; E = '  -3.49%  ' }
    @{ Row = 49; D = '7.900'; E = '  -1.35%  ' }
    @{ Row = 50; D = '0.05150'; E = '  -1.59%  ' }
    @{ Row = 51; D = '0.4207'; E = '  -1.17%  ' }
)

$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

$dataRange.ClearFormats()
